# M08 Froze Encoder 123
# Update Epoch Accuracy values (recomputed after re-running training) and
# refresh the stale DisplayOutputs object repr (new Python object memory
# address) in column A for the later block of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Epoch Accuracy (column B) values for epochs 2-50 (rows 4-52)
$ws.Range("B4").Value = 0.34375
$ws.Range("B5").Value = 0.34375
$ws.Range("B6").Value = 0.34375
$ws.Range("B7").Value = 0.359375
$ws.Range("B8").Value = 0.328125
$ws.Range("B9").Value = 0.34375
$ws.Range("B10").Value = 0.34375
$ws.Range("B11").Value = 0.34375
$ws.Range("B14").Value = 0.375
$ws.Range("B15").Value = 0.4375
$ws.Range("B16").Value = 0.28125
$ws.Range("B17").Value = 0.25
$ws.Range("B18").Value = 0.3125
$ws.Range("B19").Value = 0.265625
$ws.Range("B20").Value = 0.25
$ws.Range("B21").Value = 0.265625
$ws.Range("B22").Value = 0.25
$ws.Range("B23").Value = 0.265625
$ws.Range("B24").Value = 0.265625
$ws.Range("B25").Value = 0.265625
$ws.Range("B26").Value = 0.265625
$ws.Range("B27").Value = 0.265625
$ws.Range("B28").Value = 0.265625
$ws.Range("B29").Value = 0.265625
$ws.Range("B30").Value = 0.265625
$ws.Range("B31").Value = 0.265625
$ws.Range("B32").Value = 0.265625
$ws.Range("B33").Value = 0.265625
$ws.Range("B34").Value = 0.265625
$ws.Range("B35").Value = 0.265625
$ws.Range("B36").Value = 0.28125
$ws.Range("B37").Value = 0.28125
$ws.Range("B38").Value = 0.265625
$ws.Range("B39").Value = 0.265625
$ws.Range("B40").Value = 0.265625
$ws.Range("B41").Value = 0.265625
$ws.Range("B42").Value = 0.265625
$ws.Range("B43").Value = 0.265625
$ws.Range("B44").Value = 0.265625
$ws.Range("B45").Value = 0.265625
$ws.Range("B46").Value = 0.265625
$ws.Range("B47").Value = 0.265625
$ws.Range("B48").Value = 0.265625
$ws.Range("B49").Value = 0.265625
$ws.Range("B50").Value = 0.265625
$ws.Range("B51").Value = 0.265625
$ws.Range("B52").Value = 0.265625

# Update DisplayOutputs object memory address text (column A) and accuracy
# values (column B) for rows 102-118
$ws.Range("A102").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("A103").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B103").Value = 0.34375
$ws.Range("A104").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B104").Value = 0.203125
$ws.Range("A105").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B105").Value = 0.28125
$ws.Range("A106").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B106").Value = 0.265625
$ws.Range("A107").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B107").Value = 0.140625
$ws.Range("A108").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B108").Value = 0.28125
$ws.Range("A109").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B109").Value = 0.125
$ws.Range("A110").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B110").Value = 0.1875
$ws.Range("A111").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("A112").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("A113").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B113").Value = 0.25
$ws.Range("A114").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B114").Value = 0.125
$ws.Range("A115").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B115").Value = 0.125
$ws.Range("A116").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B116").Value = 0.15625
$ws.Range("A117").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B117").Value = 0.203125
$ws.Range("A118").Value = "<__main__.DisplayOutputs object at 0x7f4b706763a0>"
$ws.Range("B118").Value = 0.2459016393442623
